$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins / Losses / Ties) in row 1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, bordered, centered) from an existing
# header cell onto the new header cells so they match the rest of row 1
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record for every player row (2-48):
# every row gets the team's season record -> 86 wins, 76 losses, 0 ties
$ws.Range("AD2:AD48").Value = 86
$ws.Range("AE2:AE48").Value = 76
$ws.Range("AF2:AF48").Value = 0
